$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Replace the "Null" text hyperlink-less cell D4 with an actual hyperlink
$url = "https://github.com/justinvandelaar/InteractiveList"
$ws.Range("D4").Value = $url
$ws.Hyperlinks.Add($ws.Range("D4"), $url, "", "", $url) | Out-Null
$ws.Range("D4").Style = "Hyperlink"

# Update selection to match the committed state
$ws.Range("B15").Select() | Out-Null
